$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - new user "Allam"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Allam"
$ws.Range("C5").Value = "allam@gmail.com"
$ws.Range("D5").Value = "pbkdf2:sha256:260000`$jhqElP07ppJlhh3v`$8e6af83a284632ae86d6b2c5f430b6b0efd7e235db92dbca4e9b7793331b61c8"
$ws.Range("E5").Value = "2025-06-26T21:36:13.811742"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = $true

# Row 6 - new user "Allam1"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Allam1"
$ws.Range("C6").Value = "allam100@gmail.com"
$ws.Range("D6").Value = "Allam123@"
$ws.Range("E6").Value = "2025-06-27T02:08:49.212447"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = $true
